$d = $word.ActiveDocument

$d.Content.Find.Execute("Descentralización del trabajo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Descentralización del trabajo", 2)

$d.Content.Find.Execute("Asignación de autoridad", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Asignación de autoridad", 2)

$d.Content.Find.Execute("5 acciones que vincula la organización con el proceso administrativo", $true, $false, $false, $false, $false,
                         $true, 1, $false, "5 acciones que vincula la organización con el proceso administrativo", 2)

$d.Content.Find.Execute("trabajadores reciban en promedio al menos dos veces la compensación mínima establecida por ley en todas sus ubicaciones. Adicionalmente, se otorgan bonos anuales por desempeño y compensaciones variables a colaboradores de tiempo completo, medio tiempo y temporal.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "trabajadores reciban en promedio al menos dos veces la compensación mínima establecida por ley en todas sus ubicaciones. Adicionalmente, se otorgan bonos anuales por desempeño y compensaciones variables a colaboradores de tiempo completo, medio tiempo y temporal.", 2)
